{"js": "// Insert a new \"Meta description\" paragraph right after the title (Heading1),\n// then, at the end of the document, remove the duplicated title paragraph and\n// turn the duplicated meta-description paragraph into the new image \"Prompt\".\n\nconst TITLE_TEXT = \"Play Dragon's Luck Power Reels for Free - Review\";\nconst META_LABEL = \"Meta description\";\nconst META_REST =\n  \": Discover Dragon's Luck Power Reels, a high-volatility slot game with excellent graphics and immersive gameplay. Play for free and win big!\";\nconst PROMPT_TEXT =\n  \"Prompt: Create a feature image for Dragon's Luck Power Reels in cartoon style featuring a happy Maya warrior with glasses. The image should incorporate the theme of Chinese tradition and the figure of the dragon in a visually appealing manner. It should also showcase the 10 reels and 30 paylines of the game, with the highest-valued symbol, the coin with the number 138, prominently displayed. The image should be optimized for use on both desktop computers and all iOS and Android mobile devices, and should capture the high volatility and distinctive gameplay features of the slot machine.\";\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"text\");\nawait context.sync();\n\n// --- 1) Insert the new \"Meta description\" paragraph right after the first\n//        occurrence of the title text (the Heading1 at the very top). ---\nconst titlePara = paragraphs.items.find((p) => p.text === TITLE_TEXT);\n\nconst metaPara = titlePara.insertParagraph(\"\", Word.InsertLocation.after);\nmetaPara.style = \"Normal\";\nawait context.sync();\n\nconst labelRange = metaPara.insertText(META_LABEL, Word.InsertLocation.start);\nlabelRange.font.bold = true;\nawait context.sync();\n\nconst restRange = metaPara.insertText(META_REST, Word.InsertLocation.end);\nrestRange.font.bold = false;\nawait context.sync();\n\n// --- 2) At the end of the document: drop the duplicated bold title paragraph\n//        and replace the duplicated italic meta-description paragraph's text\n//        with the new image-generation prompt (keeping its italic formatting). ---\nconst endParagraphs = body.paragraphs;\nendParagraphs.load(\"text\");\nawait context.sync();\n\nconst items = endParagraphs.items;\nlet trailingTitlePara = null;\nlet trailingMetaPara = null;\nfor (let i = items.length - 1; i >= 0; i--) {\n  if (trailingMetaPara === null && items[i].text === META_REST.slice(2)) {\n    trailingMetaPara = items[i];\n    continue;\n  }\n  if (trailingTitlePara === null && items[i].text === TITLE_TEXT) {\n    trailingTitlePara = items[i];\n    break;\n  }\n}\n\ntrailingMetaPara.insertText(PROMPT_TEXT, Word.InsertLocation.replace);\ntrailingTitlePara.delete();\nawait context.sync();\n", "ps1": "# Insert a new \"Meta description\" paragraph right after the title (Heading1),\n# then, at the end of the document, remove the duplicated title paragraph and\n# turn the duplicated meta-description paragraph into the new image \"Prompt\".\n\n$d = $word.ActiveDocument\n\n$titleText = \"Play Dragon's Luck Power Reels for Free - Review\"\n$metaLabel = \"Meta description\"\n$metaBody = \"Discover Dragon's Luck Power Reels, a high-volatility slot game with excellent graphics and immersive gameplay. Play for free and win big!\"\n$promptText = \"Prompt: Create a feature image for Dragon's Luck Power Reels in cartoon style featuring a happy Maya warrior with glasses. The image should incorporate the theme of Chinese tradition and the figure of the dragon in a visually appealing manner. It should also showcase the 10 reels and 30 paylines of the game, with the highest-valued symbol, the coin with the number 138, prominently displayed. The image should be optimized for use on both desktop computers and all iOS and Android mobile devices, and should capture the high volatility and distinctive gameplay features of the slot machine.\"\n\n# --- 1) Insert the new \"Meta description\" paragraph right after the first\n#        occurrence of the title text (the Heading1 at the very top). ---\n$count = $d.Paragraphs.Count\n$titleParaIndex = -1\nfor ($i = 1; $i -le $count; $i++) {\n    if ($d.Paragraphs.Item($i).Range.Text.TrimEnd() -eq $titleText) {\n        $titleParaIndex = $i\n        break\n    }\n}\n\n$titlePara = $d.Paragraphs.Item($titleParaIndex)\n$titlePara.Range.InsertParagraphAfter()\n\n$metaPara = $d.Paragraphs.Item($titleParaIndex + 1)\n$metaPara.Style = \"Normal\"\n\n$fullMetaText = $metaLabel + \": \" + $metaBody\n$metaRange = $metaPara.Range\n$metaRange.Text = $fullMetaText\n\n$boldRange = $d.Range($metaRange.Start, $metaRange.Start + $metaLabel.Length)\n$boldRange.Bold = 1\n\n# --- 2) At the end of the document: drop the duplicated bold title paragraph\n#        and replace the duplicated italic meta-description paragraph's text\n#        with the new image-generation prompt (keeping its italic formatting). ---\n$count = $d.Paragraphs.Count\n$trailingMetaIndex = -1\nfor ($i = $count; $i -ge 1; $i--) {\n    if ($d.Paragraphs.Item($i).Range.Text.TrimEnd() -eq $metaBody) {\n        $trailingMetaIndex = $i\n        break\n    }\n}\n\n$trailingTitleIndex = -1\nfor ($i = $trailingMetaIndex - 1; $i -ge 1; $i--) {\n    if ($d.Paragraphs.Item($i).Range.Text.TrimEnd() -eq $titleText) {\n        $trailingTitleIndex = $i\n        break\n    }\n}\n\n$trailingMetaPara = $d.Paragraphs.Item($trailingMetaIndex)\n$lastRange = $d.Range($trailingMetaPara.Range.Start, $trailingMetaPara.Range.End)\n$lastRange.Text = $promptText\n\n$trailingTitlePara = $d.Paragraphs.Item($trailingTitleIndex)\n$trailingTitlePara.Range.Delete()\n"}
